$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.070.09"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "1.888.04"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.66"
$ws.Range("E5").Value = "  -2.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4603"
$ws.Range("E7").Value = "  -2.99%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4071"
$ws.Range("E8").Value = "  +0.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.76"
$ws.Range("E9").Value = "  -0.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07987"
$ws.Range("E10").Value = "  -2.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9928"
$ws.Range("E11").Value = "  -3.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.72"
$ws.Range("E12").Value = "  -3.25%  "

$ws.Range("D13").Value = "1.837.94"
$ws.Range("E13").Value = "  -2.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.918"
$ws.Range("E14").Value = "  -2.85%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.076"
$ws.Range("E15").Value = "  -4.25%  "

$ws.Range("E16").Value = "  -0.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.57"
$ws.Range("E17").Value = "  -3.35%  "

$ws.Range("E18").Value = "  -2.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06559"
$ws.Range("E19").Value = "  -1.12%  "

$ws.Range("E20").Value = "  -2.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").Value = "29.076.11"
$ws.Range("E22").Value = "  +0.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.428"
$ws.Range("E23").Value = "  -2.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.48"
$ws.Range("E24").Value = "  +2.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.206"
$ws.Range("E25").Value = "  -2.90%  "

$ws.Range("D26").Value = "2.080.59"
$ws.Range("E26").Value = "  -2.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.75"
$ws.Range("E27").Value = "  -2.41%  "

$ws.Range("E28").Value = "  -2.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.102"
$ws.Range("E29").Value = "  -3.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.470"
$ws.Range("E30").Value = "  -1.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.75"
$ws.Range("E31").Value = "  -2.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.005"
$ws.Range("E32").Value = "  -0.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09330"
$ws.Range("E33").Value = "  -2.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.601"
$ws.Range("E34").Value = "  -1.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.411"
$ws.Range("E35").Value = "  -1.40%  "

$ws.Range("E36").Value = "  -2.71%  "

$ws.Range("E37").Value = "  -2.30%  "

$ws.Range("E38").Value = "  -2.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.277"
$ws.Range("E39").Value = "  -4.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.174"
$ws.Range("E40").Value = "  -2.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9987"
$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5790"
$ws.Range("E42").Value = "  -4.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1824"
$ws.Range("E43").Value = "  -4.11%  "

$ws.Range("E44").Value = "  -4.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.259"
$ws.Range("E45").Value = "  -2.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07497"
$ws.Range("E46").Value = "  +3.19%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.268"
$ws.Range("E47").Value = "  +5.33%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.06"
$ws.Range("E48").Value = "  -2.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5459"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.901"
$ws.Range("E50").Value = "  -4.11%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.59"
$ws.Range("E51").Value = "  +11.93%  "

